# Update alignment in PowerPoint generation script
# - Remove the 4 trailing error ("nan" / 貨號輸入錯啦...) rows from the table
# - Increase the height of the remaining (header + 5 data) rows from 36pt (457200 EMU)
#   to 60pt (762000 EMU)
# - Re-format numeric-looking cell text to drop the trailing ".0"
# - Refresh the "輸入時間" timestamp text for each data row

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tbl = $shape.Table

# Drop the 4 trailing "nan" / error rows (originally rows 7-10).
$tbl.Rows.Item(7).Delete()
$tbl.Rows.Item(7).Delete()
$tbl.Rows.Item(7).Delete()
$tbl.Rows.Item(7).Delete()

# Grow every remaining row (header + 5 data rows) from 36pt to 60pt.
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $tbl.Rows.Item($i).Height = 60
}

$newTime = "2024-11-22 15:31:25.160000"

# Row 2: item 27893
$tbl.Cell(2,1).Shape.TextFrame.TextRange.Text = "27893"
$tbl.Cell(2,3).Shape.TextFrame.TextRange.Text = "7"
$tbl.Cell(2,7).Shape.TextFrame.TextRange.Text = $newTime

# Row 3: item 10941
$tbl.Cell(3,1).Shape.TextFrame.TextRange.Text = "10941"
$tbl.Cell(3,3).Shape.TextFrame.TextRange.Text = "12"
$tbl.Cell(3,7).Shape.TextFrame.TextRange.Text = $newTime

# Row 4: item 19103
$tbl.Cell(4,1).Shape.TextFrame.TextRange.Text = "19103"
$tbl.Cell(4,3).Shape.TextFrame.TextRange.Text = "20"
$tbl.Cell(4,7).Shape.TextFrame.TextRange.Text = $newTime

# Row 5: item 27899
$tbl.Cell(5,1).Shape.TextFrame.TextRange.Text = "27899"
$tbl.Cell(5,3).Shape.TextFrame.TextRange.Text = "10"
$tbl.Cell(5,7).Shape.TextFrame.TextRange.Text = $newTime

# Row 6: item 10642
$tbl.Cell(6,1).Shape.TextFrame.TextRange.Text = "10642"
$tbl.Cell(6,3).Shape.TextFrame.TextRange.Text = "10"
$tbl.Cell(6,7).Shape.TextFrame.TextRange.Text = $newTime
